$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replace the data in rows 2-5 (columns A..AH) with the new dataset ---
$data = @(
    ,@(45120.50694444445,22.1,15.162,4.209,46.427,38.398,17.391,57.345,26.759,11.337,17.452,18.439,19.302,5.553,17.294,24.297,14.517,3.762,2.458,255.893,48.146,15.963,31.894,16.657,2.022,28.586,14.1,12.689,14.848,19.074,3.64,50.56,8.856,19.957)
    ,@(45120.51388888889,20.178,14.467,1.944,43.239,35.765,15.88,61.805,24.432,10.655,16.046,17.406,18.213,5.073,15.79,22.348,13.415,1.65,1.176,233.037,44.186,14.575,29.463,15.683,1.891,29.869,12.874,11.559,13.553,18.166,1.294,55.879,8.132,18.223)
    ,@(45120.52083333334,7.687,5.32,1.032,16.296,13.482,6.05,28.025,9.308,3.994,5.976,6.636,6.837,1.937,6.015,8.474,5.269,1.007,0.556,84.262,17.01,5.552,11.195,6.051,0.648,12.896,4.904,4.497,5.256,6.881,0.784,25.472,3.025,6.943)
    ,@(45120.52777777778,16.81,12.29,1.12,36.26,29.98,13.23,49.51,20.36,8.98,13.44,14.65,15.34,4.23,13.16,18.67,11.14,0.83,0.72,192.96,36.73,12.15,24.61,13.14,1.61,24.14,10.73,9.56,11.23,15.33,0.5600000000000001,44.58,6.81,15.19)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $rowIndex = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowVals[$c]
    }
}

# --- 2. Delete row 6 (dataset shrank from 5 data rows to 4 data rows) ---
$ws.Rows.Item(6).Delete()

# --- 3. Update column widths (character widths that round-trip to the target stored width) ---
$widthCols = @(2,3,5,6,7,8,9,10,11,12,13,15,16,17,21,22,23,24,26,27,28,29,30,32,34)
foreach ($col in $widthCols) {
    $ws.Columns.Item($col).ColumnWidth = 7.14
}
$ws.Columns.Item(20).ColumnWidth = 8.14
